$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 13
$ws.Cells.Item(2, 3).Value = 23999
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0.386484375
$ws.Cells.Item(2, 6).Value = 0.3895451294573561
$ws.Cells.Item(2, 7).Value = -0.019375
$ws.Cells.Item(2, 8).Value = 14.45002427304607
$ws.Cells.Item(2, 9).Value = 0.1604804911245413
$ws.Cells.Item(2, 10).Value = 36.47288402410519
$ws.Cells.Item(2, 11).Value = @"
[[ 3.34323263e-09  1.37253246e-10  4.48108361e-10  6.15178112e-08
  -1.30800621e-07  2.85545815e-07]
 [ 5.02919562e-12  8.07071605e-09 -5.58822965e-10 -6.15817844e-07
  -3.52227708e-09  8.81095146e-09]
 [ 4.07220739e-11 -5.58822965e-10  4.12996262e-09  1.13759582e-06
  -2.84871382e-08  4.12720119e-08]
 [ 1.00247748e-08 -6.15817844e-07  1.13759582e-06  3.58721905e-04
  -7.01264694e-06  9.74130915e-06]
 [-2.24539386e-08 -3.52227708e-09 -2.84871382e-08 -7.01264694e-06
   1.52451431e-05 -3.33815694e-05]
 [ 4.91751047e-08  8.81095146e-09  4.12720119e-08  9.74130915e-06
  -3.33815694e-05  7.32290460e-05]]
"@
$ws.Cells.Item(2, 12).Value = -0.1995191749448844
$ws.Cells.Item(2, 13).Value = 26.3703380250261
$ws.Cells.Item(2, 14).Value = 2.20283536651915
$ws.Cells.Item(2, 15).Value = 0.02382671508122759
$ws.Cells.Item(2, 16).Value = 0.006749365298061104
$ws.Cells.Item(2, 17).Value = 0.3310399441423864
$ws.Cells.Item(2, 18).Value = 0.01955543805103895
$ws.Cells.Item(2, 19).Value = -0.01269619445519747
$ws.Cells.Item(2, 20).Value = @"
[[ 1.24920316e-08 -2.66258375e-09 -2.22615906e-09 -1.49386697e-07
  -4.04619988e-07  5.55845343e-07]
 [-8.21472967e-11  4.71065846e-08  3.36866008e-08  3.19540080e-06
   1.04204789e-08 -1.48795569e-07]
 [-8.18142685e-11  3.36866008e-08  4.35307493e-08  5.08036906e-06
   1.03785922e-08 -1.80423383e-07]
 [-8.38119236e-09  3.19540080e-06  5.08036906e-06  6.29240770e-04
   1.06320688e-06 -2.04914025e-05]
 [-6.44692578e-08  1.04204789e-08  1.03785922e-08  1.06320688e-06
   7.33224650e-06 -9.97490996e-06]
 [ 8.76656103e-08 -1.48795569e-07 -1.80423383e-07 -2.04914025e-05
  -9.97490996e-06  1.42190943e-05]]
"@
$ws.Cells.Item(2, 21).Value = -0.5971245776998035
$ws.Cells.Item(2, 22).Value = 12.31190331064625
$ws.Cells.Item(2, 23).Value = 1.354351398701363
$ws.Cells.Item(2, 24).Value = 0.05103341984294531
$ws.Cells.Item(2, 25).Value = 0.01199494571473045
$ws.Cells.Item(2, 26).Value = 0.6582492635724109
$ws.Cells.Item(2, 27).Value = -0.003382710706803671
$ws.Cells.Item(2, 28).Value = 0.02401624282889409
$ws.Cells.Item(2, 29).Value = -0.5029742135069977
$ws.Cells.Item(2, 30).Value = 1.008995566122331
$ws.Cells.Item(2, 31).Value = 1.513508827211741
$ws.Cells.Item(2, 32).Value = 0
$ws.Cells.Item(2, 33).Value = 90.00002493749997
$ws.Cells.Item(2, 34).Value = 25.5029
$ws.Cells.Item(2, 35).Value = 6.775275000000001
$ws.Cells.Item(2, 36).Value = 11.09699999999999
$ws.Cells.Item(2, 37).Value = 0.7
$ws.Cells.Item(2, 38).Value = 50177
$ws.Cells.Item(2, 39).Value = 1
$ws.Cells.Item(2, 40).Value = $true
$ws.Cells.Item(2, 41).Value = $true

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 13
$ws.Cells.Item(3, 3).Value = 23999
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0.392421875
$ws.Cells.Item(3, 6).Value = 0.400292173376275
$ws.Cells.Item(3, 7).Value = -0.01953125
$ws.Cells.Item(3, 8).Value = 14.62449438552168
$ws.Cells.Item(3, 9).Value = 0.1623516657921827
$ws.Cells.Item(3, 10).Value = 36.84504553115146
$ws.Cells.Item(3, 11).Value = @"
[[ 3.35175926e-09  1.28001579e-10  4.15417175e-10  5.58318835e-08
  -1.34719743e-07  2.94366234e-07]
 [ 5.26209957e-12  8.88806639e-09 -5.87098428e-10 -6.72231035e-07
  -3.72901214e-09  9.22736706e-09]
 [ 4.09373824e-11 -5.87098428e-10  4.63376893e-09  1.25435286e-06
  -2.89774425e-08  3.96452922e-08]
 [ 9.89045436e-09 -6.72231035e-07  1.25435286e-06  3.90010225e-04
  -7.00076293e-06  9.10870927e-06]
 [-2.51814660e-08 -3.72901214e-09 -2.89774425e-08 -7.00076293e-06
   1.72730742e-05 -3.78453047e-05]
 [ 5.51818967e-08  9.22736706e-09  3.96452922e-08  9.10870927e-06
  -3.78453047e-05  8.30719805e-05]]
"@
$ws.Cells.Item(3, 12).Value = -0.2010163679237675
$ws.Cells.Item(3, 13).Value = 26.40176067418819
$ws.Cells.Item(3, 14).Value = 2.202868413645728
$ws.Cells.Item(3, 15).Value = 0.02379835718048294
$ws.Cells.Item(3, 16).Value = 0.006886945032165119
$ws.Cells.Item(3, 17).Value = 0.3468727376230543
$ws.Cells.Item(3, 18).Value = 0.01733133263297782
$ws.Cells.Item(3, 19).Value = -0.007712573263520741
$ws.Cells.Item(3, 20).Value = @"
[[ 1.47588238e-08 -3.34991796e-09 -2.86757041e-09 -1.94658624e-07
  -4.83347956e-07  6.64837188e-07]
 [-2.35046738e-10  6.11723244e-08  4.48093074e-08  4.03960244e-06
   3.06636715e-08 -2.19347764e-07]
 [-2.36129767e-10  4.48093074e-08  5.69809328e-08  6.30338174e-06
   3.08051347e-08 -2.59061552e-07]
 [-2.29813545e-08  4.03960244e-06  6.30338174e-06  7.42697023e-04
   2.99810205e-06 -2.76652204e-05]
 [-7.87418636e-08  3.06636715e-08  3.08051347e-08  2.99810205e-06
   9.13851492e-06 -1.25024966e-05]
 [ 1.07605980e-07 -2.19347764e-07 -2.59061552e-07 -2.76652204e-05
  -1.25024966e-05  1.79486112e-05]]
"@
$ws.Cells.Item(3, 21).Value = -0.590791683442049
$ws.Cells.Item(3, 22).Value = 12.31944797701489
$ws.Cells.Item(3, 23).Value = 1.354265119328042
$ws.Cells.Item(3, 24).Value = 0.05100216599723046
$ws.Cells.Item(3, 25).Value = 0.01255269357825202
$ws.Cells.Item(3, 26).Value = 0.7116730859819469
$ws.Cells.Item(3, 27).Value = -0.00827104435114088
$ws.Cells.Item(3, 28).Value = 0.02919019357618715
$ws.Cells.Item(3, 29).Value = -0.5029745943664593
$ws.Cells.Item(3, 30).Value = 1.008999058819861
$ws.Cells.Item(3, 31).Value = 1.513504243010262
$ws.Cells.Item(3, 32).Value = 0
$ws.Cells.Item(3, 33).Value = 90.00002493749997
$ws.Cells.Item(3, 34).Value = 25.5029
$ws.Cells.Item(3, 35).Value = 6.775275000000001
$ws.Cells.Item(3, 36).Value = 11.09699999999999
$ws.Cells.Item(3, 37).Value = 0.7
$ws.Cells.Item(3, 38).Value = 49718
$ws.Cells.Item(3, 39).Value = 1
$ws.Cells.Item(3, 40).Value = $true
$ws.Cells.Item(3, 41).Value = $true

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 13
$ws.Cells.Item(4, 3).Value = 23999
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 0.387890625
$ws.Cells.Item(4, 6).Value = 0.4110392172951939
$ws.Cells.Item(4, 7).Value = -0.0196875
$ws.Cells.Item(4, 8).Value = 14.56170173348895
$ws.Cells.Item(4, 9).Value = 0.1619134804583391
$ws.Cells.Item(4, 10).Value = 36.69081052818206
$ws.Cells.Item(4, 11).Value = @"
[[ 3.51032102e-09  1.43850620e-10  4.72628467e-10  6.46940778e-08
  -1.37853718e-07  3.00942025e-07]
 [ 5.11455338e-12  8.13095588e-09 -5.64205883e-10 -6.18028927e-07
  -3.57378407e-09  8.92016113e-09]
 [ 4.11333504e-11 -5.64205883e-10  4.22034083e-09  1.15430810e-06
  -2.87090429e-08  4.13131437e-08]
 [ 1.00545011e-08 -6.18028927e-07  1.15430810e-06  3.61122524e-04
  -7.01736094e-06  9.67339440e-06]
 [-2.27828182e-08 -3.57378407e-09 -2.87090429e-08 -7.01736094e-06
   1.54332824e-05 -3.37943453e-05]
 [ 4.98965054e-08  8.92016113e-09  4.13131437e-08  9.67339440e-06
  -3.37943453e-05  7.41371991e-05]]
"@
$ws.Cells.Item(4, 12).Value = -0.1991434820987736
$ws.Cells.Item(4, 13).Value = 26.28438158632956
$ws.Cells.Item(4, 14).Value = 2.202827074110476
$ws.Cells.Item(4, 15).Value = 0.02390463434166339
$ws.Cells.Item(4, 16).Value = 0.006790984717259054
$ws.Cells.Item(4, 17).Value = 0.3297938198398468
$ws.Cells.Item(4, 18).Value = 0.01876500467619497
$ws.Cells.Item(4, 19).Value = -0.01077562721389261
$ws.Cells.Item(4, 20).Value = @"
[[ 1.50703138e-08 -3.57707002e-09 -3.02325755e-09 -2.04756072e-07
  -4.90023195e-07  6.74368875e-07]
 [-7.79237288e-11  3.56761610e-08  2.59276469e-08  2.36959281e-06
   1.00315821e-08 -1.16477559e-07]
 [-7.75214240e-11  2.59276469e-08  3.27242747e-08  3.66484272e-06
   9.98005127e-09 -1.38282453e-07]
 [-7.63672493e-09  2.36959281e-06  3.66484272e-06  4.36826879e-04
   9.83150926e-07 -1.50494823e-05]
 [-4.56639124e-08  1.00315821e-08  9.98005127e-09  9.83150926e-07
   5.23342209e-06 -7.12963840e-06]
 [ 6.21695776e-08 -1.16477559e-07 -1.38282453e-07 -1.50494823e-05
  -7.12963840e-06  1.01977769e-05]]
"@
$ws.Cells.Item(4, 21).Value = -0.5909360910168804
$ws.Cells.Item(4, 22).Value = 12.30151158501554
$ws.Cells.Item(4, 23).Value = 1.354267086614866
$ws.Cells.Item(4, 24).Value = 0.05107653042275818
$ws.Cells.Item(4, 25).Value = 0.01240245114469744
$ws.Cells.Item(4, 26).Value = 0.7047853649007698
$ws.Cells.Item(4, 27).Value = -0.002483259171899697
$ws.Cells.Item(4, 28).Value = 0.02167290368083439
$ws.Cells.Item(4, 29).Value = -0.5029748664092858
$ws.Cells.Item(4, 30).Value = 1.008996548441067
$ws.Cells.Item(4, 31).Value = 1.513503260684987
$ws.Cells.Item(4, 32).Value = 0
$ws.Cells.Item(4, 33).Value = 90.00002493749997
$ws.Cells.Item(4, 34).Value = 25.5029
$ws.Cells.Item(4, 35).Value = 6.775275000000001
$ws.Cells.Item(4, 36).Value = 11.09699999999999
$ws.Cells.Item(4, 37).Value = 0.7
$ws.Cells.Item(4, 38).Value = 49867
$ws.Cells.Item(4, 39).Value = 1
$ws.Cells.Item(4, 40).Value = $true
$ws.Cells.Item(4, 41).Value = $true

# Apply the same style as the header row's formatted cells (bold, centered, thin border)
# to the new index column A2:A4, matching style index 1 used by row 1.
$ws.Range("B1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# Reset auto-fit row heights (writing multi-line matrix text can trigger
# an explicit row height); restore default row sizing to match a plain
# data import with no manual row-height overrides.
$ws.Range("A2:A4").EntireRow.AutoFit()
